$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "85.871.20"
$ws.Range("E2").Value = "  +6.77%  "

# Row 3
$ws.Range("D3").Value = "3.314.54"
$ws.Range("E3").Value = "  +3.15%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").Value = "'220.00"
$ws.Range("E5").Value = "  +4.64%  "

# Row 6
$ws.Range("D6").Value = "'635.46"
$ws.Range("E6").Value = "  +0.44%  "

# Row 7
$ws.Range("D7").Value = "'0.322"
$ws.Range("E7").Value = "  +17.06%  "

# Row 8
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "  -2.25%  "

# Row 10
$ws.Range("D10").Value = "3.316.90"
$ws.Range("E10").Value = "  +3.15%  "

# Row 11
$ws.Range("D11").Value = "'0.595"
$ws.Range("E11").Value = "  -2.95%  "

# Row 12
$ws.Range("D12").Value = "'0.0000274"
$ws.Range("E12").Value = "  +4.47%  "

# Row 13
$ws.Range("E13").Value = "  +0.13%  "

# Row 14
$ws.Range("D14").Value = "3.933.27"
$ws.Range("E14").Value = "  +3.00%  "

# Row 15
$ws.Range("D15").Value = "'34.08"
$ws.Range("E15").Value = "  +4.77%  "

# Row 16
$ws.Range("D16").Value = "'5.40"
$ws.Range("E16").Value = "  -0.24%  "

# Row 17
$ws.Range("D17").Value = "85.674.82"
$ws.Range("E17").Value = "  +6.24%  "

# Row 18
$ws.Range("D18").Value = "3.314.52"
$ws.Range("E18").Value = "  +2.61%  "

# Row 19
$ws.Range("D19").Value = "'14.58"
$ws.Range("E19").Value = "  -0.05%  "

# Row 20
$ws.Range("D20").Value = "'3.17"
$ws.Range("E20").Value = "  +5.80%  "

# Row 21
$ws.Range("D21").Value = "'441.21"
$ws.Range("E21").Value = "  -1.73%  "

# Row 22
$ws.Range("D22").Value = "'9.12"
$ws.Range("E22").Value = "  -2.61%  "

# Row 23
$ws.Range("D23").Value = "'5.23"
$ws.Range("E23").Value = "  -2.04%  "

# Row 24
$ws.Range("D24").Value = "'7.30"
$ws.Range("E24").Value = "  +5.51%  "

# Row 25
$ws.Range("D25").Value = "'5.44"
$ws.Range("E25").Value = "  +13.72%  "

# Row 26
$ws.Range("D26").Value = "'12.22"
$ws.Range("E26").Value = "  +11.19%  "

# Row 27
$ws.Range("D27").Value = "3.490.10"
$ws.Range("E27").Value = "  +2.68%  "

# Row 28
$ws.Range("D28").Value = "'78.26"
$ws.Range("E28").Value = "  +0.52%  "

# Row 29
$ws.Range("D29").Value = "'0.0000130"
$ws.Range("E29").Value = "  +4.19%  "

# Row 30
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.03%  "

# Row 31
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.170"
$ws.Range("E31").Value = "  +35.11%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'609.85"
$ws.Range("E32").Value = "  +8.68%  "

# Row 33
$ws.Range("D33").Value = "'9.23"
$ws.Range("E33").Value = "  -0.42%  "

# Row 34
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = "  +2.03%  "

# Row 36
$ws.Range("D36").Value = "'2.04"
$ws.Range("E36").Value = "  +0.11%  "

# Row 37
$ws.Range("D37").Value = "'0.151"
$ws.Range("E37").Value = "  -2.71%  "

# Row 38
$ws.Range("D38").Value = "'23.30"
$ws.Range("E38").Value = "  -2.37%  "

# Row 39
$ws.Range("D39").Value = "'6.42"
$ws.Range("E39").Value = "  +11.01%  "

# Row 40
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.32%  "

# Row 41
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.417"
$ws.Range("E41").Value = "  -0.54%  "

# Row 42
$ws.Range("D42").Value = "'21.26"
$ws.Range("E42").Value = "  +4.40%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.06"
$ws.Range("E43").Value = "  +11.64%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'3.06"
$ws.Range("E44").Value = "  +12.12%  "

# Row 45
$ws.Range("D45").Value = "'158.64"
$ws.Range("E45").Value = "  -4.49%  "

# Row 46
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("D47").Value = "'188.67"
$ws.Range("E47").Value = "  -1.81%  "

# Row 48
$ws.Range("E48").Value = "  +1.80%  "

# Row 49
$ws.Range("D49").Value = "'45.19"
$ws.Range("E49").Value = "  +2.97%  "

# Row 50
$ws.Range("D50").Value = "'0.792"
$ws.Range("E50").Value = "  -1.42%  "

# Row 51
$ws.Range("D51").Value = "'26.37"
$ws.Range("E51").Value = "  +3.14%  "
